# ------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计"
#    containing the per-fund holding detail for the new quarter.
# 2. Prepend a new summary row for "2022-Q1" at the top of the "总计"
#    (totals) worksheet's data, pushing the existing rows down by one
#    and renumbering the index column.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: create the new "2022-Q1" worksheet positioned right after
# "2021-Q4" (and therefore right before "总计").
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("2021-Q4"))
$newSheet.Name = "2022-Q1"

# NOTE: worksheet object references captured before a sheet is
# inserted can become stale once the workbook's sheet collection
# changes, so every worksheet we need is (re)fetched by name right
# before it is used.

# Copy the look & feel (fonts / borders / alignment) of the header row
# and the index column from the "2021-Q4" sheet so the new sheet's
# styling matches the rest of the workbook.
$wb.Worksheets.Item("2021-Q4").Range("B1:H1").Copy()
$wb.Worksheets.Item("2022-Q1").Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Worksheets.Item("2021-Q4").Range("A2").Copy()
$wb.Worksheets.Item("2022-Q1").Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet = $wb.Worksheets.Item("2022-Q1")

# ---- header row -----------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# ---- data rows --------------------------------------------------------
# Each row: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @("202023", "南方优选成长混合A", "41.34", "61.77", "1.21", "0.5002", 10),
    @("004818", "国寿安保目标策略灵活配置混合A", "4.06", "36.45", "3.99", "0.1620", 1),
    @("001672", "国寿安保智慧生活股票", "3.56", "85.91", "4.22", "0.1502", 1),
    @("519677", "银河定投宝中证腾讯济安价值100A股指数", "2.74", "91.56", "1.21", "0.0332", 10),
    @("005206", "南方优选成长混合C", "2.40", "61.77", "1.21", "0.0290", 10),
    @("004819", "国寿安保目标策略灵活配置混合C", "0.57", "36.45", "3.99", "0.0227", 1),
    @("001731", "广发百发大数据策略价值灵活配置混合A", "0.24", "88.87", "4.07", "0.0098", 2),
    @("001732", "广发百发大数据策略价值灵活配置混合E", "0.24", "88.87", "4.07", "0.0098", 2)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    # A column: numeric running index (0-based), style already copied above
    $newSheet.Cells.Item($r, 1).Value = $i

    # B column: fund code - must stay text so leading zeros survive
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    # C column: fund name (plain text, never numeric-looking)
    $newSheet.Cells.Item($r, 3).Value = $row[1]

    # D-G columns: numeric-looking figures stored as text in the source
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    # H column: rank, numeric
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ------------------------------------------------------------------
# Step 2: insert the "2022-Q1" summary row at the top of "总计"'s data
# (row 2), pushing the existing rows down and renumbering column A.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Shift the existing 3 data rows (currently rows 2-4) down to rows 3-5,
# working bottom-up so we never overwrite data we still need to read.
for ($r = 4; $r -ge 2; $r--) {
    $dst = $r + 1
    $totalSheet.Cells.Item($dst, 2).Value = $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($dst, 3).Value = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($dst, 4).Value = $totalSheet.Cells.Item($r, 4).Value2
}

# Make sure the new row 5 (A5) has the same styling as the other index
# cells (A2:A4) before we write its value.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Renumber the index column (A2:A5) to 0,1,2,3
for ($r = 2; $r -le 5; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Write the brand-new "2022-Q1" summary values into row 2.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 0.92
